$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new review-log entry as row 7 (mirrors rows 2-6: A=flag, B=name,
# C=time opened, D=time finished review, E=file folder path).
$ws.Range("A7").Value = 0
$ws.Range("B7").Value = "Alex"
$ws.Range("C7").Value = 45063.85536992636
$ws.Range("D7").Value = 45063.85543055303
$ws.Range("E7").Value = "C:/Users/Alex/Documents/Builds/ReviewApp/Test"

# Carry over the same look-and-feel used by the previous log row (bold /
# centered / bordered flag cell, date-time number format on the two
# timestamp cells) by copying row 6's formatting onto the new row.
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)
